$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.261.90'
$ws.Range("E2").Value = '  -0.78%  '
$ws.Range("D3").Value = '1.903.67'
$ws.Range("E3").Value = '  +1.18%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("B5").Value = 'XRP'
$ws.Range("C5").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.693'
$ws.Range("E5").Value = '  +9.57%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '245.62'
$ws.Range("E6").Value = '  +0.82%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '41.66'
$ws.Range("E8").Value = '  -3.00%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.350'
$ws.Range("E9").Value = '  +5.06%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '52.97'
$ws.Range("E10").Value = '  +12.19%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0726'
$ws.Range("E11").Value = '  +3.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0996'
$ws.Range("E12").Value = '  +0.21%  '
$ws.Range("D13").Value = '2.178.89'
$ws.Range("E13").Value = '  +1.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '12.33'
$ws.Range("E14").Value = '  +0.14%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.715'
$ws.Range("E15").Value = '  +4.28%  '
$ws.Range("D16").Value = '1.903.34'
$ws.Range("E16").Value = '  +1.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.84'
$ws.Range("E17").Value = '  +1.30%  '
$ws.Range("D18").Value = '35.274.28'
$ws.Range("E18").Value = '  -0.82%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.39'
$ws.Range("E19").Value = '  +1.20%  '
$ws.Range("E20").Value = '  +1.84%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '241.12'
$ws.Range("E21").Value = '  -1.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '12.62'
$ws.Range("E22").Value = '  +1.56%  '
$ws.Range("E23").Value = '  -0.21%  '
$ws.Range("E24").Value = '  -0.30%  '
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.28'
$ws.Range("E26").Value = '  +11.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '169.17'
$ws.Range("E27").Value = '  -1.02%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.58'
$ws.Range("E28").Value = '  +2.84%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.40'
$ws.Range("E29").Value = '  +2.93%  '
$ws.Range("E30").Value = '  +4.35%  '
$ws.Range("E32").Value = '  +2.44%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.960'
$ws.Range("E33").Value = '  +2.12%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0572'
$ws.Range("E34").Value = '  +1.29%  '
$ws.Range("E35").Value = '  -0.30%  '
$ws.Range("E36").Value = '  +0.56%  '
$ws.Range("E37").Value = '  -0.54%  '
$ws.Range("E38").Value = '  -1.52%  '
$ws.Range("E39").Value = '  -1.24%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0688'
$ws.Range("E40").Value = '  +15.57%  '
$ws.Range("E41").Value = '  -0.36%  '
$ws.Range("E42").Value = '  +2.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '16.18'
$ws.Range("E43").Value = '  +5.38%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '90.38'
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").Value = '1.345.37'
$ws.Range("E45").Value = '  -0.80%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.45'
$ws.Range("E46").Value = '  +4.74%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.85'
$ws.Range("E47").Value = '  +3.25%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '12.58'
$ws.Range("E48").Value = '  -4.59%  '
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("E50").Value = '  +1.95%  '
$ws.Range("E51").Value = '  -1.92%  '
